$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 399) holds the "Förändrad" date serial value.
# It changes from 45205 (2023-10-06) to 45206 (2023-10-07) for every row.
$ws.Range("C2:C399").Value = 45206
